# Generate Report for Handback
# Overview sheet: re-sort rows by file name and refresh the handback status
# for 89df6b3e-7a68-49f8-8273-30b129d815f6.md (now "Handed back: in sync
# with en-US" instead of "Ready for handoff"), which also moves it ahead of
# cabe8da7-b78e-46d2-8df4-48ff53588842.md in row order.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "89df6b3e-7a68-49f8-8273-30b129d815f6.md"
$ov.Range("B2").Value = "e2e\89df6b3e-7a68-49f8-8273-30b129d815f6.md"
$ov.Range("G2").Value = "2016-08-26 22:47:02"

$ov.Range("A3").Value = "cabe8da7-b78e-46d2-8df4-48ff53588842.md"
$ov.Range("B3").Value = "e2e\cabe8da7-b78e-46d2-8df4-48ff53588842.md"
$ov.Range("E3").Value = "Handed back: in sync with en-US"
$ov.Range("F3").Value = "Handed back: in sync with en-US"
$ov.Range("G3").Value = "2016-08-26 22:45:51"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "89df6b3e-7a68-49f8-8273-30b129d815f6.md"
$zh.Range("G2").Value = "89df6b3e-7a68-49f8-8273-30b129d815f6.69cd22beefe8b94dee6ffc7afc6df686b01d2cc0.zh-cn.xlf"
$zh.Range("H2").Value = "2016-08-26 22:46:55"
$zh.Range("I2").Value = "89df6b3e-7a68-49f8-8273-30b129d815f6.md"
$zh.Range("J2").Value = "89df6b3e-7a68-49f8-8273-30b129d815f6.69cd22beefe8b94dee6ffc7afc6df686b01d2cc0.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-26 22:47:19"

$zh.Range("A3").Value = "cabe8da7-b78e-46d2-8df4-48ff53588842.md"
$zh.Range("C3").Value = "Handed back: in sync with en-US"
$zh.Range("G3").Value = "cabe8da7-b78e-46d2-8df4-48ff53588842.2e1843878eaeb4e45cb9723c90d56ef6ac4a171a.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-26 22:45:46"
$zh.Range("I3").Value = "cabe8da7-b78e-46d2-8df4-48ff53588842.md"
$zh.Range("J3").Value = "cabe8da7-b78e-46d2-8df4-48ff53588842.2e1843878eaeb4e45cb9723c90d56ef6ac4a171a.zh-cn.xlf"
$zh.Range("P3").Value = ""

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "89df6b3e-7a68-49f8-8273-30b129d815f6.md"
$de.Range("G2").Value = "89df6b3e-7a68-49f8-8273-30b129d815f6.69cd22beefe8b94dee6ffc7afc6df686b01d2cc0.de-de.xlf"
$de.Range("H2").Value = "2016-08-26 22:47:02"
$de.Range("I2").Value = "89df6b3e-7a68-49f8-8273-30b129d815f6.md"
$de.Range("J2").Value = "89df6b3e-7a68-49f8-8273-30b129d815f6.69cd22beefe8b94dee6ffc7afc6df686b01d2cc0.de-de.xlf"
$de.Range("K2").Value = "2016-08-26 22:47:26"

$de.Range("A3").Value = "cabe8da7-b78e-46d2-8df4-48ff53588842.md"
$de.Range("C3").Value = "Handed back: in sync with en-US"
$de.Range("G3").Value = "cabe8da7-b78e-46d2-8df4-48ff53588842.2e1843878eaeb4e45cb9723c90d56ef6ac4a171a.de-de.xlf"
$de.Range("H3").Value = "2016-08-26 22:45:51"
$de.Range("I3").Value = "cabe8da7-b78e-46d2-8df4-48ff53588842.md"
$de.Range("J3").Value = "cabe8da7-b78e-46d2-8df4-48ff53588842.2e1843878eaeb4e45cb9723c90d56ef6ac4a171a.de-de.xlf"
$de.Range("P3").Value = ""
